$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (deepseek_r1_1_5b, Group_A)
$ws.Range("J2").Value = 104

# Row 3 (deepseek_r1_1_5b, Group_B)
$ws.Range("J3").Value = 447
$ws.Range("L3").Value = 67
$ws.Range("M3").Value = 0

# Row 4 (deepseek_r1_1_5b, Group_C)
$ws.Range("J4").Value = 400
$ws.Range("L4").Value = 38
$ws.Range("M4").Value = 0

# Row 5 (deepseek_r1_8b, Group_A)
$ws.Range("J5").Value = 101

# Row 6 (deepseek_r1_8b, Group_B)
$ws.Range("D6").Value = 26
$ws.Range("J6").Value = 80
$ws.Range("L6").Value = 163

# Row 7 (deepseek_r1_14b, Group_A)
$ws.Range("J7").Value = 146

# Row 8 (deepseek_r1_14b, Group_B)
$ws.Range("D8").Value = 3
$ws.Range("J8").Value = 9

# Row 9 (deepseek_r1_14b, Group_C)
$ws.Range("D9").Value = 1
$ws.Range("J9").Value = 3

# Row 10 (deepseek_r1_32b, Group_A)
$ws.Range("J10").Value = 109
